$d = $word.ActiveDocument

# Texts of the five paragraphs (inside the "Soutien a la politique de formation"
# table) that lose their justified (w:jc w:val="both") alignment.
$targets = @(
    "Nombre de stagiaires",
    "Affaires Maritimes",
    "LPM/ENSM",
    "Étrangers",
    "Total de jours/stagiaires"
)

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    foreach ($needle in $targets) {
        if ($t.StartsWith($needle)) {
            $p.Alignment = 0
            break
        }
    }
}

# Remove the two manual line breaks after "Total de jours/stagiaires : /".
$d.Content.Find.Execute(
    "Total de jours/stagiaires : /" + [char]11 + [char]11,
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Total de jours/stagiaires : /", 2)

# Merge the now-empty trailing paragraph in that table cell into the
# "Total de jours/stagiaires" paragraph (it used to be a leftover blank
# paragraph after the two manual line breaks).
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Total de jours/stagiaires")) {
        $d.Paragraphs.Item($i + 1).Range.Delete()
        break
    }
}
